# Export-Evaluationsmatrix.xlsx — fill in the three previously-empty
# evaluation criteria rows (3-5) with their labels and scores.
#
# NOTE on string insertion order: the shared-strings table lists new
# entries in first-seen order, so we write the A5 label before A3/A4 to
# reproduce "Darstellungsmöglichkeiten", "Einfache Handhabung", "Kosten"
# (in that sequence) at the tail of sharedStrings.xml.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 label first (see note above), then rows 3 and 4.
$ws.Range("A5").Value = "Darstellungsmöglichkeiten"
$ws.Range("A3").Value = "Einfache Handhabung"
$ws.Range("A4").Value = "Kosten"

# Row 3: Einfache Handhabung — Gewicht 2, HTML 10, PDF 10
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 10
$ws.Range("D3").Value = 10

# Row 4: Kosten — Gewicht 1, HTML 10, PDF 10
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 10
$ws.Range("D4").Value = 10

# Row 5: Darstellungsmöglichkeiten — Gewicht 1, HTML 9, PDF 8
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 9
$ws.Range("D5").Value = 8

# Match the author's final selection (active cell B5).
$ws.Range("B5").Select()
